# Auto-generated Excel COM-interop script to apply profit-column value updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1625.8823
$ws.Range("I19").Value = 3210
$ws.Range("J19").Value = 517
$ws.Range("K19").Value = 3210
$ws.Range("L19").Value = 517
$ws.Range("M19").Value = -3035
$ws.Range("N19").Value = -867
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H80").Value = 9369885
$ws.Range("I80").Value = 642.5333000000001
$ws.Range("J80").Value = 22146126
$ws.Range("K80").Value = 1927.5999
$ws.Range("L80").Value = 66438378
$ws.Range("M80").Value = -929.5999000000002
$ws.Range("N80").Value = -66440374
$ws.Range("H83").Value = 9369885
$ws.Range("I83").Value = 642.5333000000001
$ws.Range("J83").Value = 22146126
$ws.Range("K83").Value = 5782.7997
$ws.Range("L83").Value = 199315134
$ws.Range("M83").Value = -790.7997000000005
$ws.Range("N83").Value = -199325118
$ws.Range("H116").Value = 3714.6155
$ws.Range("I116").Value = 979
$ws.Range("K116").Value = 979
$ws.Range("M116").Value = 2463
$ws.Range("H137").Value = 31870.031
$ws.Range("I137").Value = 1506.1538
$ws.Range("K137").Value = 4518.4614
$ws.Range("M137").Value = -1968.4614

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2181.4524
$ws.Range("I45").Value = 2030.258
$ws.Range("K45").Value = 2030.258
$ws.Range("M45").Value = -1653.258
$ws.Range("H110").Value = 703.6667
$ws.Range("I110").Value = 711
$ws.Range("K110").Value = 711
$ws.Range("M110").Value = 1334
$ws.Range("H132").Value = 11887.68
$ws.Range("I132").Value = 1695.2703
$ws.Range("J132").Value = 40896.848
$ws.Range("K132").Value = 5085.810899999999
$ws.Range("L132").Value = 122690.544
$ws.Range("M132").Value = -2555.810899999999
$ws.Range("N132").Value = -127750.544

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 924.7059
$ws.Range("I20").Value = 987.2
$ws.Range("J20").Value = 835.4286
$ws.Range("K20").Value = 987.2
$ws.Range("L20").Value = 835.4286
$ws.Range("M20").Value = -740.2
$ws.Range("N20").Value = -1329.4286
$ws.Range("H107").Value = 1823.625
$ws.Range("I107").Value = 1755.5714
$ws.Range("K107").Value = 1755.5714
$ws.Range("M107").Value = 164.4286
$ws.Range("H134").Value = 55903.3
$ws.Range("I134").Value = 91575.164
$ws.Range("K134").Value = 274725.492
$ws.Range("M134").Value = -272190.492

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12993.237
$ws.Range("J31").Value = 3230.2273
$ws.Range("L31").Value = 3230.2273
$ws.Range("N31").Value = -3820.2273
$ws.Range("H34").Value = 12993.237
$ws.Range("J34").Value = 3230.2273
$ws.Range("L34").Value = 3230.2273
$ws.Range("N34").Value = -3634.2273
$ws.Range("H99").Value = 16132698
$ws.Range("I99").Value = 3425.6
$ws.Range("J99").Value = 45458650
$ws.Range("K99").Value = 3425.6
$ws.Range("L99").Value = 45458650
$ws.Range("M99").Value = -1927.6
$ws.Range("N99").Value = -45461646
$ws.Range("H107").Value = 1009.68
$ws.Range("J107").Value = 1179.2142
$ws.Range("L107").Value = 1179.2142
$ws.Range("N107").Value = -5019.2142
$ws.Range("H126").Value = 16132698
$ws.Range("I126").Value = 3425.6
$ws.Range("J126").Value = 45458650
$ws.Range("K126").Value = 10276.8
$ws.Range("L126").Value = 136375950
$ws.Range("M126").Value = -7806.799999999999
$ws.Range("N126").Value = -136380890

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2500159.8
$ws.Range("I4").Value = 191.6
$ws.Range("K4").Value = 574.8
$ws.Range("M4").Value = -462.8
$ws.Range("H40").Value = 193.72728
$ws.Range("I40").Value = 97.888885
$ws.Range("J40").Value = 625
$ws.Range("K40").Value = 391.55554
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -322.55554
$ws.Range("N40").Value = -2638
$ws.Range("H80").Value = 13895.75
$ws.Range("J80").Value = 21313.4
$ws.Range("L80").Value = 63940.2
$ws.Range("N80").Value = -65812.20000000001
$ws.Range("H81").Value = 5507.3
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5507.3
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 16521.9
$ws.Range("M81").Value = ""
$ws.Range("N81").Value = -18767.9
$ws.Range("H83").Value = 13895.75
$ws.Range("J83").Value = 21313.4
$ws.Range("L83").Value = 191820.6
$ws.Range("N83").Value = -201180.6
$ws.Range("H84").Value = 5507.3
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5507.3
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 49565.7
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = -60797.7
$ws.Range("H86").Value = 41666920
$ws.Range("I86").Value = 258.66666
$ws.Range("J86").Value = 83333580
$ws.Range("K86").Value = 775.9999799999999
$ws.Range("L86").Value = 250000740
$ws.Range("M86").Value = 410.0000200000001
$ws.Range("N86").Value = -250003112
$ws.Range("H89").Value = 41666920
$ws.Range("I89").Value = 258.66666
$ws.Range("J89").Value = 83333580
$ws.Range("K89").Value = 2327.99994
$ws.Range("L89").Value = 750002220
$ws.Range("M89").Value = 3600.00006
$ws.Range("N89").Value = -750014076
$ws.Range("H92").Value = 41666960
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H134").Value = 2532.1365
$ws.Range("I134").Value = 1305.2
$ws.Range("K134").Value = 3915.6
$ws.Range("M134").Value = 1154.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6556250
$ws.Range("J12").Value = 6000000
$ws.Range("L12").Value = 6000000
$ws.Range("N12").Value = -6000280
$ws.Range("H107").Value = 8547320
$ws.Range("J107").Value = 38461690
$ws.Range("L107").Value = 38461690
$ws.Range("N107").Value = -38465530
$ws.Range("H122").Value = 666666700
$ws.Range("I122").Value = 333333340
$ws.Range("J122").Value = 1000000000
$ws.Range("K122").Value = 1000000020
$ws.Range("L122").Value = 3000000000
$ws.Range("M122").Value = -999997570
$ws.Range("N122").Value = -3000004900
$ws.Range("H132").Value = 47223.91
$ws.Range("I132").Value = 51318.43
$ws.Range("J132").Value = 40609.69
$ws.Range("K132").Value = 153955.29
$ws.Range("L132").Value = 121829.07
$ws.Range("M132").Value = -151425.29
$ws.Range("N132").Value = -126889.07

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2866
$ws.Range("I82").Value = 2450
$ws.Range("J82").Value = 3490
$ws.Range("K82").Value = 2450
$ws.Range("L82").Value = 3490
$ws.Range("M82").Value = -2089
$ws.Range("N82").Value = -4212
$ws.Range("H85").Value = 2866
$ws.Range("I85").Value = 2450
$ws.Range("J85").Value = 3490
$ws.Range("K85").Value = 2450
$ws.Range("L85").Value = 3490
$ws.Range("M85").Value = -1202
$ws.Range("N85").Value = -5986
$ws.Range("H100").Value = 2312
$ws.Range("I100").Value = 1374.75
$ws.Range("J100").Value = 3249.25
$ws.Range("K100").Value = 1374.75
$ws.Range("L100").Value = 3249.25
$ws.Range("M100").Value = -833.75
$ws.Range("N100").Value = -4331.25
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H137").Value = 35400.715
$ws.Range("J137").Value = 35400.715
$ws.Range("L137").Value = 35400.715
$ws.Range("N137").Value = -45600.715
$ws.Range("H141").Value = 47950
$ws.Range("J141").Value = 47950
$ws.Range("L141").Value = 47950
$ws.Range("N141").Value = -58310

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3497146.8
$ws.Range("I107").Value = 676.5
$ws.Range("K107").Value = 2029.5
$ws.Range("M107").Value = -109.5
$ws.Range("H132").Value = 2143
$ws.Range("I132").Value = 1526.4445
$ws.Range("K132").Value = 4579.333500000001
$ws.Range("M132").Value = -2049.333500000001

Write-Output "Applied all profit sheet updates."